# Update benchmark: 2025-10-24 06:36:27 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Clear D/F (and related) cells that now hold no value
$ws.Range("D3").Value = ""
$ws.Range("F3").Value = ""

$ws.Range("D4").Value = ""
$ws.Range("F4").Value = ""

$ws.Range("D5").Value = ""
$ws.Range("F5").Value = ""

$ws.Range("D6").Value = ""

$ws.Range("D8").Value = ""
$ws.Range("F8").Value = ""

$ws.Range("D9").Value = ""
$ws.Range("F9").Value = ""

$ws.Range("D10").Value = ""
$ws.Range("F10").Value = ""

$ws.Range("D11").Value = ""

$ws.Range("D12").Value = ""

$ws.Range("D13").Value = ""
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("F13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"

$ws.Range("D14").Value = ""
$ws.Range("F14").Value = ""
